# Apply targeted cell-text updates to the single-column results table.
# Each row of the table holds one metric value in its sole cell; we
# address rows/cells directly (rather than a blind Find/Replace) because
# several values (e.g. "100", "0", "0.00004") repeat across rows.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $newText) {
    $cell = $table.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $newText
}

# Simple scalar value fix-ups
Set-CellText $t 1  "0M"
Set-CellText $t 2  "0M"
Set-CellText $t 3  "0M"
Set-CellText $t 4  "121"
Set-CellText $t 5  "0.00002"
Set-CellText $t 9  "0.00005"
Set-CellText $t 12 "0.00484"

# Rows 44-46 previously held a full tab-delimited detail line crammed into
# a single run; collapse each back down to just its leading summary value.
Set-CellText $t 44 "100"
Set-CellText $t 45 "0"
Set-CellText $t 46 "106"
